# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price refresh values to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 20002096
$ws.Range("I100").Value = 1300
$ws.Range("J100").Value = 25002294
$ws.Range("K100").Value = 1300
$ws.Range("L100").Value = 25002294
$ws.Range("M100").Value = -759
$ws.Range("N100").Value = -25003376

$ws.Range("H132").Value = 3369652.2
$ws.Range("I132").Value = 2876.2856
$ws.Range("J132").Value = 22223598
$ws.Range("K132").Value = 8628.856800000001
$ws.Range("L132").Value = 66670794
$ws.Range("M132").Value = -6098.856800000001
$ws.Range("N132").Value = -66675854

$ws.Range("H138").Value = 3164.5789
$ws.Range("I138").Value = 2926.3809
$ws.Range("J138").Value = 3458.8235
$ws.Range("K138").Value = 8779.1427
$ws.Range("L138").Value = 10376.4705
$ws.Range("M138").Value = -3639.1427
$ws.Range("N138").Value = -20656.4705

$ws.Range("H141").Value = 2749.2222
$ws.Range("I141").Value = 2032.4
$ws.Range("J141").Value = 6333.3335
$ws.Range("K141").Value = 6097.200000000001
$ws.Range("L141").Value = 19000.0005
$ws.Range("M141").Value = -917.2000000000007
$ws.Range("N141").Value = -29360.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3466177.8
$ws.Range("I61").Value = 1667668.5
$ws.Range("J61").Value = 14706860
$ws.Range("K61").Value = 1667668.5
$ws.Range("L61").Value = 14706860
$ws.Range("M61").Value = -1667456.5
$ws.Range("N61").Value = -14707284

$ws.Range("H102").Value = 6632.227
$ws.Range("I102").Value = 2206.75
$ws.Range("J102").Value = 18433.5
$ws.Range("K102").Value = 2206.75
$ws.Range("L102").Value = 18433.5
$ws.Range("M102").Value = -584.75
$ws.Range("N102").Value = -21677.5

$ws.Range("H132").Value = 30784194
$ws.Range("I132").Value = 40874536
$ws.Range("J132").Value = 5558335.5
$ws.Range("K132").Value = 122623608
$ws.Range("L132").Value = 16675006.5
$ws.Range("M132").Value = -122621078

$ws.Range("H136").Value = 3466177.8
$ws.Range("I136").Value = 1667668.5
$ws.Range("J136").Value = 14706860
$ws.Range("K136").Value = 5003005.5
$ws.Range("L136").Value = 44120580
$ws.Range("M136").Value = -5000455.5
$ws.Range("N136").Value = -44125680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1913.89
$ws.Range("I86").Value = 1934.7368
$ws.Range("J86").Value = 1517.8
$ws.Range("K86").Value = 1934.7368
$ws.Range("L86").Value = 1517.8
$ws.Range("M86").Value = -811.7367999999999
$ws.Range("N86").Value = -3763.8

$ws.Range("H89").Value = 1913.89
$ws.Range("I89").Value = 1934.7368
$ws.Range("J89").Value = 1517.8
$ws.Range("K89").Value = 9673.683999999999
$ws.Range("L89").Value = 7589
$ws.Range("M89").Value = -4057.683999999999
$ws.Range("N89").Value = -18821

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2475689.2
$ws.Range("I31").Value = 3206815.2
$ws.Range("J31").Value = 1899650.5
$ws.Range("K31").Value = 3206815.2
$ws.Range("L31").Value = 1899650.5
$ws.Range("M31").Value = -3206520.2
$ws.Range("N31").Value = -1900240.5

$ws.Range("H34").Value = 2475689.2
$ws.Range("I34").Value = 3206815.2
$ws.Range("J34").Value = 1899650.5
$ws.Range("K34").Value = 3206815.2
$ws.Range("L34").Value = 1899650.5
$ws.Range("M34").Value = -3206613.2
$ws.Range("N34").Value = -1900054.5

$ws.Range("H43").Value = 40000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 40000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 40000
$ws.Range("N43").Value = -40368

$ws.Range("H58").Value = 1981390.5
$ws.Range("I58").Value = 6696.9414
$ws.Range("J58").Value = 7576355.5
$ws.Range("K58").Value = 6696.9414
$ws.Range("L58").Value = 7576355.5
$ws.Range("M58").Value = -6493.9414
$ws.Range("N58").Value = -7576761.5

$ws.Range("H62").Value = 8335860
$ws.Range("I62").Value = 2282.7646
$ws.Range("J62").Value = 19233616
$ws.Range("K62").Value = 2282.7646
$ws.Range("L62").Value = 19233616
$ws.Range("M62").Value = -1658.7646
$ws.Range("N62").Value = -19234864

$ws.Range("H65").Value = 8335860
$ws.Range("I65").Value = 2282.7646
$ws.Range("J65").Value = 19233616
$ws.Range("K65").Value = 11413.823
$ws.Range("L65").Value = 96168080
$ws.Range("M65").Value = -8293.823
$ws.Range("N65").Value = -96174320

$ws.Range("H101").Value = 40000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 40000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490

$ws.Range("H136").Value = 1981390.5
$ws.Range("I136").Value = 6696.9414
$ws.Range("J136").Value = 7576355.5
$ws.Range("K136").Value = 20090.8242
$ws.Range("L136").Value = 22729066.5
$ws.Range("M136").Value = -17540.8242
$ws.Range("N136").Value = -22734166.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = $null

$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = $null

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = $null

$ws.Range("H121").Value = 2783118.5
$ws.Range("I121").Value = 1973.6
$ws.Range("J121").Value = 3514998.5
$ws.Range("K121").Value = 5920.799999999999
$ws.Range("L121").Value = 10544995.5
$ws.Range("M121").Value = -4610.799999999999
$ws.Range("N121").Value = -10547615.5

$ws.Range("H129").Value = 2108.8215
$ws.Range("I129").Value = 1304.8
$ws.Range("J129").Value = 2555.5
$ws.Range("K129").Value = 3914.4
$ws.Range("L129").Value = 7666.5
$ws.Range("M129").Value = 1085.6
$ws.Range("N129").Value = -17666.5

$ws.Range("H131").Value = 969.1053000000001
$ws.Range("I131").Value = 487.14285
$ws.Range("J131").Value = 1018
$ws.Range("K131").Value = 1461.42855
$ws.Range("L131").Value = 3054
$ws.Range("M131").Value = 3578.57145
$ws.Range("N131").Value = -13134

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12702.632
$ws.Range("I80").Value = 5831.25
$ws.Range("J80").Value = 17700
$ws.Range("K80").Value = 5831.25
$ws.Range("L80").Value = 17700
$ws.Range("M80").Value = -4833.25

$ws.Range("H83").Value = 12702.632
$ws.Range("I83").Value = 5831.25
$ws.Range("J83").Value = 17700
$ws.Range("K83").Value = 29156.25
$ws.Range("L83").Value = 88500
$ws.Range("M83").Value = -24164.25

$ws.Range("H97").Value = 10417464
$ws.Range("I97").Value = 796.86664
$ws.Range("J97").Value = 27778576
$ws.Range("K97").Value = 796.86664
$ws.Range("L97").Value = 27778576
$ws.Range("M97").Value = -300.86664
$ws.Range("N97").Value = -27779568

$ws.Range("H132").Value = 11675570
$ws.Range("I132").Value = 13034677
$ws.Range("J132").Value = 9093267
$ws.Range("K132").Value = 39104031
$ws.Range("L132").Value = 27279801
$ws.Range("M132").Value = -39101501
$ws.Range("N132").Value = -27284861

$ws.Range("H140").Value = 45000
$ws.Range("I140").Value = 30000
$ws.Range("J140").Value = 60000
$ws.Range("K140").Value = 30000
$ws.Range("L140").Value = 60000
$ws.Range("M140").Value = -24820
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2990
$ws.Range("I100").Value = 3140
$ws.Range("J100").Value = 2840
$ws.Range("K100").Value = 3140
$ws.Range("L100").Value = 2840
$ws.Range("M100").Value = -2599
$ws.Range("N100").Value = -3922

$ws.Range("H127").Value = 36118.332
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 36118.332
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 36118.332
$ws.Range("N127").Value = -46038.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4167
$ws.Range("I136").Value = 1400.4615
$ws.Range("J136").Value = 5194.5713
$ws.Range("K136").Value = 4201.3845
$ws.Range("L136").Value = 15583.7139
$ws.Range("M136").Value = -1651.3845
$ws.Range("N136").Value = -20683.7139

Write-Host "Applied scheduled-runner price refresh to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
